$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update "Förändrad" (changed) date column C for rows 2-15 from 45185 (2023-09-16)
# to 45204 (2023-10-05), preserving existing date formatting/style.
$ws.Range("C2:C15").Value = 45204
